# Generate Report for Handback
# Refresh the "Latest Handback DateTime" for the 650f1eed-... localization
# item on both the zh-cn and de-de report sheets (column K, row 2) to the
# timestamp of the newly generated handback report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("K2").Value = "2016-10-13 13:56:06"
$wsDeDe.Range("K2").Value = "2016-10-13 13:56:23"
